# The document originally contains three reviewer comments:
#   - two comments anchored on the "EDUCATION & CERTIFICATIONS" heading
#   - one comment anchored on the word "pediatric"
# This change removes all of them (and their associated comment-range /
# comment-reference markup), leaving the underlying text untouched.

$d = $word.ActiveDocument

# Delete every comment in the document. Iterate backwards since the
# collection re-indexes as items are removed.
for ($i = $d.Comments.Count; $i -ge 1; $i--) {
    $d.Comments.Item($i).Delete()
}
